# Updated symbol list on Sun Dec 18 06:43:16 UTC 2022 with GitHub Actions
#
# Source data cells in column D look numeric ("246.43", "0.02985", ...) but
# are stored as TEXT in the workbook (no leading "=", general display of
# trailing zeros that a real number would drop). Excel's COM `Value` setter
# auto-detects such strings and silently converts them to real numbers,
# which would both change the cell type and drop significant trailing
# zeros (e.g. "0.09260" -> 0.0926). Prefixing with a literal leading
# apostrophe forces Excel to keep them as text (exactly like typing
# '0.09260 into a cell), and the apostrophe itself is not stored as part
# of the value. That trick does stamp the cell with a "quote prefix" style
# though, so we immediately reset the style back to Normal to keep the
# cell's formatting identical to the original (no style changes are in
# scope for this edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNumber($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

# --- Simple price (column D) updates ---
Set-TextNumber "D2"  "246.24"
Set-TextNumber "D4"  "5.595"
Set-TextNumber "D6"  "3.406"
Set-TextNumber "D7"  "6.475"
Set-TextNumber "D9"  "1.069"
Set-TextNumber "D19" "0.006274"
Set-TextNumber "D20" "0.001053"
Set-TextNumber "D24" "3.980"
Set-TextNumber "D25" "2.120"
Set-TextNumber "D27" "0.1291"
Set-TextNumber "D40" "0.04182"
Set-TextNumber "D41" "0.007131"
Set-TextNumber "D42" "0.003505"
Set-TextNumber "D43" "0.1044"
Set-TextNumber "D44" "0.009814"
Set-TextNumber "D45" "0.00005636"
Set-TextNumber "D47" "0.6807"
Set-TextNumber "D49" "0.00002103"

# Row 48: price update + removal of "Worstin24h" suffix
Set-TextNumber "D48" "0.02858"
$ws.Range("E48").Value = "47BOLOBOLO"

# --- Rows 10-18: coin listing shifted up by one position, with new
#     price/volume data; row 18 wraps around to "One" and gains the
#     "Worstin24h" marker that previously sat on row 48 (BOLO). ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextNumber "D10" "0.1434"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextNumber "D11" "0.07424"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextNumber "D12" "0.03184"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextNumber "D13" "0.02979"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextNumber "D14" "0.09260"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextNumber "D15" "0.001668"
$ws.Range("E15").Value = "14BitForexTokenBF"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextNumber "D16" "3.274"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextNumber "D17" "0.04696"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextNumber "D18" "0.0005756"
$ws.Range("E18").Value = "17OneONEWorstin24h"
